$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number by Excel;
# mark them as Text first so the literal string (e.g. trailing zeros) is preserved.
$textCells = @("D5","D6","D7","D9","D12","D13","D14","D20","D21","D23","D24","D25","D26","D30","D32","D33","D34","D35","D37","D38","D39","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin rows in sheet order.
$ws.Range("D2").Value = '51.768.34'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '2.830.87'
$ws.Range("E3").Value = '  +2.07%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '351.72'
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").Value = '113.30'
$ws.Range("E6").Value = '  +5.02%  '

$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +1.73%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +5.86%  '

$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("E11").Value = '  -0.48%  '

$ws.Range("D12").Value = '0.0850'
$ws.Range("E12").Value = '  +1.72%  '

$ws.Range("D13").Value = '20.02'
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").Value = '7.79'
$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '3.277.84'
$ws.Range("E15").Value = '  +2.10%  '

$ws.Range("E16").Value = '  +6.39%  '

$ws.Range("D17").Value = '2.840.29'
$ws.Range("E17").Value = '  +2.57%  '

$ws.Range("D18").Value = '51.818.71'
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("E19").Value = '  +12.20%  '

$ws.Range("D20").Value = '7.61'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = '13.37'
$ws.Range("E21").Value = '  +1.86%  '

$ws.Range("D22").Value = '0.0₃0974'
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").Value = '70.57'
$ws.Range("E23").Value = '  +1.09%  '

$ws.Range("D24").Value = '268.60'
$ws.Range("E24").Value = '  +1.21%  '

$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  +1.82%  '

$ws.Range("D26").Value = '26.27'
$ws.Range("E26").Value = '  +1.04%  '

$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("E28").Value = '  +0.69%  '

$ws.Range("E29").Value = '  +4.05%  '

$ws.Range("D30").Value = '38.86'
$ws.Range("E30").Value = '  +6.35%  '

$ws.Range("E31").Value = '  +2.94%  '

$ws.Range("D32").Value = '6.34'
$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("D33").Value = '52.84'
$ws.Range("E33").Value = '  +1.82%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0898'
$ws.Range("E34").Value = '  +8.65%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '0.0453'
$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("E36").Value = '  +2.36%  '

$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").Value = '19.14'
$ws.Range("E38").Value = '  +4.40%  '

$ws.Range("D39").Value = '3.22'
$ws.Range("E39").Value = '  +2.62%  '

$ws.Range("E40").Value = '  +2.91%  '

$ws.Range("E41").Value = '  +2.11%  '

$ws.Range("D42").Value = '2.52'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D43").Value = '122.49'
$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("D44").Value = '22.31'
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("D45").Value = '2.21'
$ws.Range("E45").Value = '  +1.03%  '

$ws.Range("D46").Value = '2.190.41'
$ws.Range("E46").Value = '  +4.30%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '2.52'
$ws.Range("E47").Value = '  +9.22%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '3.51'
$ws.Range("E48").Value = '  +8.19%  '

$ws.Range("D49").Value = '0.242'
$ws.Range("E49").Value = '  +20.04%  '

$ws.Range("D50").Value = '0.956'
$ws.Range("E50").Value = '  +6.48%  '

$ws.Range("D51").Value = '5.52'
$ws.Range("E51").Value = '  +2.31%  '
